# Applies the "gh-pages output" refresh to the 苏州-漫展信息 workbook.
#
# Observed transformation (from the OOXML diff):
#   - The oldest event row (row 2: "2024-09-07 苏州·AME动漫嘉年华...") was
#     removed from both the "展览" (exhibitions) sheet and the "全部类型"
#     (all-types) sheet, which causes every following row to shift up by
#     one position.
#   - Column A (the running index, 0/1/2/3/...) is NOT touched by the
#     shift - it always stays the static sequence tied to the row number.
#   - Columns B-E, G, H, I simply carry over the content that used to sit
#     one row below.
#   - Column F ("想去人数" / interest count) is refreshed with newly
#     scraped counts for the surviving rows (most unchanged or +a few,
#     one -1).
#   - The sheets "演出" and "本地生活" are untouched.

$wb = $excel.ActiveWorkbook

# New "想去人数" (interest count) values scraped for the surviving rows,
# keyed by their FINAL row number on the "展览" sheet (rows 2..34).
$exhibitionCounts = @{
    2  = 1027;  3  = 13518; 4  = 40;   5  = 1027;  6  = 18;
    7  = 1736;  8  = 135;   9  = 121;  10 = 79;    11 = 40;
    12 = 31;    13 = 13531; 14 = 336;  15 = 601;   16 = 8957;
    17 = 8;     18 = 8042;  19 = 254;  20 = 9;     21 = 148;
    22 = 425;   23 = 146;   24 = 7;    25 = 21;    26 = 1021;
    27 = 13;    28 = 18;    29 = 394;  30 = 206;   31 = 182;
    32 = 378;   33 = 95;    34 = 5222
}

# Same counts, keyed by final row number on the "全部类型" sheet
# (rows 2..36 - two extra "演出" rows sit between what are rows 28/29 on
# the "展览" sheet, so the row numbers diverge from row 29 onward).
$allTypesCounts = @{
    2  = 1027;  3  = 13518; 4  = 40;   5  = 1027;  6  = 18;
    7  = 1736;  8  = 135;   9  = 121;  10 = 79;    11 = 40;
    12 = 31;    13 = 13531; 14 = 336;  15 = 601;   16 = 8957;
    17 = 8;     18 = 8042;  19 = 254;  20 = 9;     21 = 148;
    22 = 425;   23 = 146;   24 = 7;    25 = 21;    26 = 1021;
    27 = 13;    28 = 18;    29 = 38;   30 = 2;     31 = 394;
    32 = 206;   33 = 182;   34 = 378;  35 = 95;    36 = 5222
}

function Remove-TopEventRow {
    param(
        [string]$SheetName,
        [hashtable]$NewCounts
    )

    $ws = $wb.Worksheets.Item($SheetName)


    # Used range row count before the edit (includes the header row).
    $lastRow = $ws.UsedRange.Rows.Count

    # Shift rows 2..lastRow-1 up by one: column A (the index) stays put;
    # columns B-E, G-I are copied down from the row below; column F gets
    # the freshly scraped count.
    for ($r = 2; $r -le ($lastRow - 1); $r++) {
        $srcRow = $r + 1

        $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($srcRow, 2).Value()   # B 开始时间
        $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($srcRow, 3).Value()   # C 名称
        $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($srcRow, 4).Value()   # D 地点
        $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($srcRow, 5).Value()   # E 具体时间范围
        $ws.Cells.Item($r, 6).Value = $NewCounts[$r]                      # F 想去人数 (refreshed)
        $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($srcRow, 7).Value()  # G 最低票价
        $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($srcRow, 8).Value()  # H Link
        $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($srcRow, 9).Value()  # I Cover
    }

    # The old bottom-most data row is now a duplicate of row (lastRow-1);
    # drop it so the sheet shrinks by exactly one row.
    $ws.Rows.Item($lastRow).Delete()
}

Remove-TopEventRow "展览" $exhibitionCounts
Remove-TopEventRow "全部类型" $allTypesCounts
